$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the new "Tree Depth" table data (rows 55-60, columns C:D)
$labels = @("Toxic", "Severely Toxic", "Obscene", "Threat", "Insult", "Identity Hate")
$depths = @(8, 6, 5, 9, 5, 6)

for ($i = 0; $i -lt 6; $i++) {
    $row = 55 + $i
    $ws.Range("C$row").Value = $labels[$i]
    $ws.Range("D$row").Value = $depths[$i]
}

# Add axis titles to the Precision chart (second chart object on the sheet)
$chartObj = $ws.ChartObjects(2)
$chart = $chartObj.Chart

$catAxis = $chart.Axes(1)
$catAxis.HasTitle = $true
$catAxis.AxisTitle.Text = "Tree Depth"

$valAxis = $chart.Axes(2)
$valAxis.HasTitle = $true
$valAxis.AxisTitle.Text = "Precision"
